# Apply the "Updated symbol list" data refresh to the crypto price table.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h), F=Data, G=Hora.
# D/E/G cells hold numeric-looking text (Excel would otherwise auto-convert
# them to numbers/percentages), so each is written with a leading apostrophe
# to force a text literal and then restyled back to "Normal" so no stray
# number-format/quote-prefix style sticks to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = 'D2'; Value = "'293.71"; Numeric = $true },
    @{ Cell = 'E2'; Value = "'-5.50%"; Numeric = $true },
    @{ Cell = 'G2'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D3'; Value = "'31.36"; Numeric = $true },
    @{ Cell = 'E3'; Value = "'-3.47%"; Numeric = $true },
    @{ Cell = 'G3'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D4'; Value = "'5.090"; Numeric = $true },
    @{ Cell = 'E4'; Value = "'-4.73%"; Numeric = $true },
    @{ Cell = 'G4'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D5'; Value = "'0.07379"; Numeric = $true },
    @{ Cell = 'E5'; Value = "'-2.95%"; Numeric = $true },
    @{ Cell = 'G5'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D6'; Value = "'7.702"; Numeric = $true },
    @{ Cell = 'E6'; Value = "'-2.19%"; Numeric = $true },
    @{ Cell = 'G6'; Value = "'9"; Numeric = $true },
    @{ Cell = 'B7'; Value = 'GateToken'; Numeric = $false },
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; Numeric = $false },
    @{ Cell = 'D7'; Value = "'3.777"; Numeric = $true },
    @{ Cell = 'E7'; Value = "'0.42%"; Numeric = $true },
    @{ Cell = 'G7'; Value = "'9"; Numeric = $true },
    @{ Cell = 'B8'; Value = 'FTXToken'; Numeric = $false },
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; Numeric = $false },
    @{ Cell = 'D8'; Value = "'1.664"; Numeric = $true },
    @{ Cell = 'E8'; Value = "'3.90%"; Numeric = $true },
    @{ Cell = 'G8'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D9'; Value = "'0.9274"; Numeric = $true },
    @{ Cell = 'E9'; Value = "'0.81%"; Numeric = $true },
    @{ Cell = 'G9'; Value = "'9"; Numeric = $true },
    @{ Cell = 'E10'; Value = "'-2.14%"; Numeric = $true },
    @{ Cell = 'G10'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D11'; Value = "'0.07140"; Numeric = $true },
    @{ Cell = 'E11'; Value = "'-6.67%"; Numeric = $true },
    @{ Cell = 'G11'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D12'; Value = "'0.07918"; Numeric = $true },
    @{ Cell = 'E12'; Value = "'-4.36%"; Numeric = $true },
    @{ Cell = 'G12'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D13'; Value = "'0.02986"; Numeric = $true },
    @{ Cell = 'E13'; Value = "'-1.74%"; Numeric = $true },
    @{ Cell = 'G13'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D14'; Value = "'0.09903"; Numeric = $true },
    @{ Cell = 'E14'; Value = "'-0.10%"; Numeric = $true },
    @{ Cell = 'G14'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D15'; Value = "'0.001490"; Numeric = $true },
    @{ Cell = 'E15'; Value = "'-1.65%"; Numeric = $true },
    @{ Cell = 'G15'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D16'; Value = "'0.006212"; Numeric = $true },
    @{ Cell = 'E16'; Value = "'0.74%"; Numeric = $true },
    @{ Cell = 'G16'; Value = "'9"; Numeric = $true },
    @{ Cell = 'E17'; Value = "'-0.26%"; Numeric = $true },
    @{ Cell = 'G17'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D18'; Value = "'2.227"; Numeric = $true },
    @{ Cell = 'E18'; Value = "'-0.93%"; Numeric = $true },
    @{ Cell = 'G18'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D19'; Value = "'0.3274"; Numeric = $true },
    @{ Cell = 'E19'; Value = "'-0.52%"; Numeric = $true },
    @{ Cell = 'G19'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D20'; Value = "'0.1349"; Numeric = $true },
    @{ Cell = 'E20'; Value = "'0.99%"; Numeric = $true },
    @{ Cell = 'G20'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D21'; Value = "'4.592"; Numeric = $true },
    @{ Cell = 'E21'; Value = "'7.98%"; Numeric = $true },
    @{ Cell = 'G21'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D22'; Value = "'0.04648"; Numeric = $true },
    @{ Cell = 'E22'; Value = "'1.62%"; Numeric = $true },
    @{ Cell = 'G22'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D23'; Value = "'0.1554"; Numeric = $true },
    @{ Cell = 'E23'; Value = "'-4.44%"; Numeric = $true },
    @{ Cell = 'G23'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D24'; Value = "'0.001221"; Numeric = $true },
    @{ Cell = 'E24'; Value = "'-0.04%"; Numeric = $true },
    @{ Cell = 'G24'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D25'; Value = "'0.004417"; Numeric = $true },
    @{ Cell = 'E25'; Value = "'-2.01%"; Numeric = $true },
    @{ Cell = 'G25'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D26'; Value = "'0.0001304"; Numeric = $true },
    @{ Cell = 'E26'; Value = "'0.15%"; Numeric = $true },
    @{ Cell = 'G26'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D27'; Value = "'0.0001879"; Numeric = $true },
    @{ Cell = 'E27'; Value = "'5.87%"; Numeric = $true },
    @{ Cell = 'G27'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G28'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G29'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G30'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G31'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G32'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G33'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G34'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G35'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G36'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G37'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G38'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D39'; Value = "'0.01653"; Numeric = $true },
    @{ Cell = 'E39'; Value = "'-6.04%"; Numeric = $true },
    @{ Cell = 'G39'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D40'; Value = "'0.04390"; Numeric = $true },
    @{ Cell = 'E40'; Value = "'-5.54%"; Numeric = $true },
    @{ Cell = 'G40'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D41'; Value = "'0.007096"; Numeric = $true },
    @{ Cell = 'E41'; Value = "'-1.33%"; Numeric = $true },
    @{ Cell = 'G41'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D42'; Value = "'0.1325"; Numeric = $true },
    @{ Cell = 'E42'; Value = "'-3.56%"; Numeric = $true },
    @{ Cell = 'G42'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D43'; Value = "'0.002107"; Numeric = $true },
    @{ Cell = 'E43'; Value = "'-6.94%"; Numeric = $true },
    @{ Cell = 'G43'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D44'; Value = "'0.01233"; Numeric = $true },
    @{ Cell = 'E44'; Value = "'-14.73%"; Numeric = $true },
    @{ Cell = 'G44'; Value = "'9"; Numeric = $true },
    @{ Cell = 'D45'; Value = "'0.00005989"; Numeric = $true },
    @{ Cell = 'E45'; Value = "'-3.40%"; Numeric = $true },
    @{ Cell = 'G45'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G46'; Value = "'9"; Numeric = $true },
    @{ Cell = 'E47'; Value = "'-15.19%"; Numeric = $true },
    @{ Cell = 'G47'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G48'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G49'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G50'; Value = "'9"; Numeric = $true },
    @{ Cell = 'G51'; Value = "'9"; Numeric = $true }
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    $cell.Value = $change.Value
    if ($change.Numeric) {
        $cell.Style = "Normal"
    }
}
